$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q1").Value = "MicroplasticImages"
$ws.Range("Q2").Value = "ALGALITA_CW_3_above500_30.jpeg"
$ws.Range("Q3").Value = "B_DW_3_above500_96.jpeg"
$ws.Range("Q4").Value = "CC_CW_1_20-250_176.jpeg"

$ws.Range("Q5").Select()
